$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Ciudades")

# Update the "last updated" timestamp in the header cell (A1)
$ws.Range("A1").Value = "Datos actualizados a 20 de Marzo de 2020 a las 20:46"

# Malaga (row 11): Recuperados 407 -> 406, Muertes 17 -> 18
$ws.Range("D11").Value = 406
$ws.Range("E11").Value = 18

# Cantabria (row 35): Recuperados 133 -> 132, Muertes 1 -> 2
$ws.Range("D35").Value = 132
$ws.Range("E35").Value = 2

# Cordoba (row 36): Recuperados 125 -> 123, Muertes 0 -> 2
$ws.Range("D36").Value = 123
$ws.Range("E36").Value = 2

# Cadiz (row 39): Recuperados 103 -> 102, Muertes 0 -> 1
$ws.Range("D39").Value = 102
$ws.Range("E39").Value = 1
